# Commit: "change multiple terminal messages type"
#
# On the "opening" sheet, several terminal-message rows have their
# "type" value (column B) changed from 2 to 3: rows 2-7 and rows 37-41.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("opening")

$rows = @(2, 3, 4, 5, 6, 7, 37, 38, 39, 40, 41)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = 3
}
